$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "28.603.74"; E = "  +2.21%  " }
    @{ Row = 3; D = "1.878.26"; E = "  +1.12%  " }
    @{ Row = 4; D = "1.011"; E = "  +0.39%  " }
    @{ Row = 5; D = "316.23"; E = "  +1.77%  " }
    @{ Row = 6; D = "1.009"; E = "  +0.94%  " }
    @{ Row = 7; D = "0.5109"; E = "  +0.71%  " }
    @{ Row = 8; D = "0.3929"; E = "  +1.70%  " }
    @{ Row = 9; D = "0.08399"; E = "  +2.88%  " }
    @{ Row = 10; D = "1.115"; E = "  +0.93%  " }
    @{ Row = 11; D = "41.75"; E = "  +0.62%  " }
    @{ Row = 12; D = "6.281"; E = "  +2.08%  " }
    @{ Row = 13; D = "1.882.89"; E = "  +1.04%  " }
    @{ Row = 14; D = "20.49"; E = "  +2.45%  " }
    @{ Row = 15; D = "7.281"; E = "  +1.86%  " }
    @{ Row = 16; D = "1.010"; E = "  +0.57%  " }
    @{ Row = 17; D = "0.00001107"; E = "  +1.71%  " }
    @{ Row = 18; D = "91.59"; E = "  +1.42%  " }
    @{ Row = 19; D = "0.06727"; E = "  +1.31%  " }
    @{ Row = 20; D = "17.75"; E = "  +1.35%  " }
    @{ Row = 21; D = "1.008"; E = "  +0.86%  " }
    @{ Row = 22; D = "5.977"; E = "  +0.71%  " }
    @{ Row = 23; D = "28.624.59"; E = "  +2.10%  " }
    @{ Row = 24; D = "11.15"; E = "  +1.59%  " }
    @{ Row = 25; D = "2.247"; E = "  +0.72%  " }
    @{ Row = 26; D = "2.099.20"; E = "  +1.00%  " }
    @{ Row = 27; D = "161.39"; E = "  +1.49%  " }
    @{ Row = 28; D = "20.86"; E = "  +1.63%  " }
    @{ Row = 29; D = "2.374"; E = "  -0.62%  " }
    @{ Row = 30; D = "127.84"; E = "  +2.19%  " }
    @{ Row = 31; D = $null; E = "  +1.51%  " }
    @{ Row = 32; D = "1.055"; E = "  +2.41%  " }
    @{ Row = 33; D = "5.814"; E = "  +0.15%  " }
    @{ Row = 34; D = "3.614"; E = "  +1.15%  " }
    @{ Row = 35; D = $null; E = "  +2.61%  " }
    @{ Row = 36; D = "0.06550"; E = "  +1.05%  " }
    @{ Row = 37; D = "0.2188"; E = "  +0.75%  " }
    @{ Row = 38; D = "8.912"; E = "  -2.74%  " }
    @{ Row = 39; D = "1.267"; E = "  +3.24%  " }
    @{ Row = 40; D = "1.198"; E = "  +2.55%  " }
    @{ Row = 41; D = "0.6481"; E = "  +1.43%  " }
    @{ Row = 42; D = "5.072"; E = "  +3.29%  " }
    @{ Row = 43; D = "11.19"; E = "  +1.11%  " }
    @{ Row = 44; D = "1.008"; E = "  +0.55%  " }
    @{ Row = 45; D = "0.6079"; E = "  +1.23%  " }
    @{ Row = 46; D = "13.09"; E = "  +1.08%  " }
    @{ Row = 47; D = "3.704"; E = "  +1.59%  " }
    @{ Row = 48; D = "2.038"; E = "  +3.22%  " }
    @{ Row = 49; D = "1.220"; E = "  +1.96%  " }
    @{ Row = 50; D = "122.54"; E = "  +1.47%  " }
    @{ Row = 51; D = "1.194"; E = "  -6.16%  " }
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        $dCell = $ws.Range("D" + $u.Row)
        $origStyle = $dCell.Style
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = $origStyle
    }
    if ($u.E -ne $null) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
